# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect that the file is now ready for handoff, with a new handoff
# timestamp and a new handoff target file (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-21 10:29:51"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-21 10:29:48"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-21 10:29:51"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
